$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26 (2025-07-30) ---
$ws.Cells.Item(26, 1).Value = 45868
$ws.Cells.Item(26, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(26, 2).Value = 6362.89990234375
$ws.Cells.Item(26, 3).Value = 6396.5400390625
$ws.Cells.Item(26, 4).Value = 6336.3798828125
$ws.Cells.Item(26, 5).Value = 6381.22998046875
$ws.Cells.Item(26, 6).Value = 5375070000
$ws.Cells.Item(26, 7).Value = -0.0012494327466496
$ws.Cells.Item(26, 8).Value = -1
$ws.Cells.Item(26, 9).Value = 'pop on smile lead the way with firstever probiotic whitening strip for a healthier brighter smile honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight why palo alto network stock is sinking today here is why exxon mobil xom is a favorite among hedge fund investor is palo alto joining the ma party plus the sp hit an earnings bump eog resource eog a top crude oil stock for passive income dr youzhi tang appointed to the cwti board of director australia secondquarter inflation drop to lowest since march supporting case for rate cut gentherm beat q revenue estimate why occidental petroleum oxy is a top choice for investment in crude oil elon musk grok ai to launch texttovideo generator this october fifth third launch private credit pact with asset manager eldridge okta inc gain amid takeover speculation oktanasdaq seeking alpha unifirst declares regular cash dividend unifirst declares regular cash dividend dr youzhi tang appointed to the cwti board of director devon energy dvn a top crude oil stock among hedge fund investor harvest fund management co ltd sell share of broadcom inc nasdaqavgo marathon petroleum mpc one of the best crude oil stock according to hedge fund kinder morgan kmi a top energy stock for passive income investor price are down but buying a home is still out of reach permian resource pr among the best crude oil stock to buy now is shell plc shel a good option for investor looking to invest in crude oil stock jpmorgan chase want to make it harder for you to use venmo toll on delaware river bridge could go up in valero energy vlo remains a top choice for hedge fund investor warren buffetts berkshire hathaway now earns every second from dividend jpmorgan reaffirms nvidia nvda a a top ai networking play capstone strategic guide janusea toward purposedriven growth through acquisition by kinective chatgpt agent pass captcha test expose flaw in bot detection system rudolph libbe group continues it steady growth ai is taking google search to new height not killing it nike is selling the most comfortable running shoe for and shopper walk in them for mile daily csp associate acquires majority stake in tamarack defense leader in market data and analytics change to federal student loan leave aspiring medical student scrambling to cover cost inventhelp inventor develops backpackstyle window cleaning device tro capital hilton to host fifth annual dc festival of magic labor day weekend the best ai crypto buy before the next bull run behavioral segment concern overshadow universal health q earnings beat alphabet googl ai strategy drive growth but analyst still see breakup potential karrueche tran stayed by deion sander side during lifesaving surgery here more about her soaring crude oil price support gain in sugar gold is this the final dip before a generational mining stock boom analyst think this former pandemic winner could surge this year cocoa price slump a the dollar rally pecheles automotive group mark year quince a bright spot in struggling dc fashion space more than double valuation to b with second raise in month trump urge the uk to drill more oil criticizing it energy policy global startup accelerator techstars officially land in tampa bay bar k bar and dog park combo close door amid financial challenge hospital in region ranked among top in florida north korea demand u recognize it a irreversible nuclear power microsofts msft ai revolution could spark a trilliondollar surge say wedbush mena renewable energy market to reach billion by driving a regional green energy boom bcc research llc mena renewable energy market to reach billion by driving a regional green energy boom bcc research llc stock market today dow fall point sp and nasdaq edge lower ahead of fed meeting and big tech earnings u and china close to extending tariff truce marketwatch brazil harvest pressure and easing tariff fear undercut coffee price google exec say employee have to be more aisavvy a competition ramp up animal shelter are closing their door to surrender a pet care cost rise how owner can afford emergency disney sued for defamation over fx series say nothing yesteryear johnson johnson nysejnj holding lowered by new advisory service llc canopy partner llc ha holding in netflix inc nasdaqnflx top lithium stock to add to your watchlist july th meeder asset management inc grows stock holding in johnson johnson nysejnj warner record name new marketing head lay off in exec realignment visualizing how mass shooting in compare with past year ivory coast president will seek a fourth term after the disqualification of rival why tapestry stock is cooling off despite coach hot streak amazon amzn price target raised to here whats driving it rosen a leading law firm encourages reddit inc investor to secure counsel before important anthropic in talk to raise fresh capital at billion valuation thermo fisher post q rise boeing report smaller loss a ceo see progress in turnaround boeing report smaller loss a ceo see progress in turnaround tax break for charitable contribution are returning for million american tax break for charitable contribution are returning for million american tax break for charitable contribution are returning for million american tax break for charitable contribution are returning for million american boeing share sink despite upbeat ndquarter earnings nick the greek to join burbank business landscape median home price by state how much do house cost rosen a leading law firm encourages reddit inc investor to secure counsel before important deadline in security class action rddt why krispy kreme stock is plummeting today inspector visit lancashire business here how each one scored for food hygiene florida struggle to keep younger worker creating headache for employer lafayette county property transfer honolulu lawsuit against fossil fuel company lead climate change legal fight utah husband wife accused of selling narcotic out of home with child letter local impact pennsylvania legislation aim to protect homebuyers in floodprone neighborhood a qa with historian marco antonio samaniego on the new binational tijuana river agreement dc council launch public hearing on billion rfk stadium proposal letter revelation show treason of those involved in russia hoax nyc subway mystery power outage ha rush hour rider vowing to steal ride a payback bipartisan bill call for annual treasury stress test pete buttigieg say gop epstein scandal reveals breakdown in societal trust trump reacts to starmer plan to recognize palestine could make the case that youre rewarding hamas from the grca can you or your organization afford to ignore professional development gop senator float broken plan to buy off american mad about tariff walking school bus and teaching student their right how educator are preparing for back to school after ice raid riverside revisits nuisance ordinance karrueche tran stayed by deion sander side during lifesaving surgery here more about her trump admin seizes judicial thorn in their side letter decriminalize possessing small amount of cannabis lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity netanyahu is reportedly planning to annex gaza strip with trump admins backing letter america need return to tradition law and order meeder asset management inc grows stake in exxon mobil corporation nysexom top lithium stock to add to your watchlist july th visualizing how mass shooting in compare with past year maxwell lawyer issue list of demand for congressional testimony lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity what american think about israel military action in gaza according to a new gallup poll what american think about israel military action in gaza according to a new gallup poll what american think about israel military action in gaza according to a new gallup poll what american think about israel military action in gaza according to a new gallup poll wire fire cause power outage in charlestown and dedham rosen a leading law firm encourages reddit inc investor to secure counsel before important deadline in security class action rddt sydney sweeneys ad spark absurd nazi claim on good morning america legacy medias latest meltdown berger montague pc investigating claim on behalf of sable offshore corp nyse soc investor after class action filing minnesota pause payment to housing stability provider amid fraud investigation what american think about israel military action in gaza according to a new gallup poll what american think about israel military action in gaza according to a new gallup poll what american think about israel military action in gaza according to a new gallup poll israeli ambassador push back on claim that israel is preventing aid distribution in gaza trump proposes to repeal landmark finding that greenhouse gas harm the public trump proposes to repeal landmark finding that greenhouse gas harm the public trump proposes to repeal landmark finding that greenhouse gas harm the public trump proposes to repeal landmark finding that greenhouse gas harm the public trump proposes to repeal landmark finding that greenhouse gas harm the public trump proposes to repeal landmark finding that greenhouse gas harm the public la county move ahead with law requiring ice agent show face badge and agency affiliation but legal hurdle loom house oversight committee reject ghislaine maxwell bid for immunity in exchange for testimony house oversight committee reject ghislaine maxwell request for immunity in exchange for her testimony biggest source of immigrant to new hampshire biggest source of immigrant to kansa biggest source of immigrant to kansa harvest fund management co ltd raise stake in alphabet inc nasdaqgoog biggest source of immigrant to alabama biggest source of immigrant to maine ucla to pay million to settle complaint of discrimination against jewish student ucla to pay million to settle complaint of discrimination against jewish student dont question my integrity spartacus return in dramatic senate floor showdown for jersey a daca recipient accidentally drove into mexico now he being fasttracked for deportation cashflow underwriting would bring fairness to mortgage lending the independent view big beautiful bill a former cdc vaccine adviser ha word for rfk jr segwayriding pest control salesperson described a pest themselves in upper township elsewhere gov walz activates minnesota national guard over st paul cyberattack gov walz activates minnesota national guard over st paul cyberattack nantucket official accuse offshore wind developer of going into hiding since trump election nantucket official accuse offshore wind developer of going into hiding since trump election nantucket official accuse offshore wind developer of going into hiding since trump election nc u senate race set to be one of nation biggest wrong way tom zirpoli big government is back commentary guest columnist an opportunity for much needed clarification together review the only one for me is you and you for me here how to see this week double meteor shower man run for mayor with really really really unfortunate name justice department fire two senior antitrust attorney alleging insubordination democraticled state sue trump admin over law that defunds planned parenthood uk will recognize palestine a a state if ceasefire not agreed on pm say mackenzie financial corp sell share of intuitive surgical inc nasdaqisrg trump say india will likely face a to tariff yearold girl reported missing in worcester'
$ws.Cells.Item(26, 10).Value = 'pop on smile lead the way with firstever probiotic whitening strip for a healthier brighter smile honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight honolulu lawsuit against fossil fuel company lead climate change legal fight why palo alto network stock is sinking today here is why exxon mobil xom is a favorite among hedge fund investor is palo alto joining the ma party plus the sp hit an earnings bump eog resource eog a top crude oil stock for passive income dr youzhi tang appointed to the cwti board of director australia secondquarter inflation drop to lowest since march supporting case for rate cut gentherm beat q revenue estimate why occidental petroleum oxy is a top choice for investment in crude oil elon musk grok ai to launch texttovideo generator this october fifth third launch private credit pact with asset manager eldridge okta inc gain amid takeover speculation oktanasdaq seeking alpha unifirst declares regular cash dividend unifirst declares regular cash dividend dr youzhi tang appointed to the cwti board of director devon energy dvn a top crude oil stock among hedge fund investor harvest fund management co ltd sell share of broadcom inc nasdaqavgo marathon petroleum mpc one of the best crude oil stock according to hedge fund kinder morgan kmi a top energy stock for passive income investor price are down but buying a home is still out of reach permian resource pr among the best crude oil stock to buy now is shell plc shel a good option for investor looking to invest in crude oil stock jpmorgan chase want to make it harder for you to use venmo toll on delaware river bridge could go up in valero energy vlo remains a top choice for hedge fund investor warren buffetts berkshire hathaway now earns every second from dividend jpmorgan reaffirms nvidia nvda a a top ai networking play capstone strategic guide janusea toward purposedriven growth through acquisition by kinective chatgpt agent pass captcha test expose flaw in bot detection system rudolph libbe group continues it steady growth ai is taking google search to new height not killing it nike is selling the most comfortable running shoe for and shopper walk in them for mile daily csp associate acquires majority stake in tamarack defense leader in market data and analytics change to federal student loan leave aspiring medical student scrambling to cover cost inventhelp inventor develops backpackstyle window cleaning device tro capital hilton to host fifth annual dc festival of magic labor day weekend the best ai crypto buy before the next bull run behavioral segment concern overshadow universal health q earnings beat alphabet googl ai strategy drive growth but analyst still see breakup potential karrueche tran stayed by deion sander side during lifesaving surgery here more about her soaring crude oil price support gain in sugar gold is this the final dip before a generational mining stock boom analyst think this former pandemic winner could surge this year cocoa price slump a the dollar rally pecheles automotive group mark year quince a bright spot in struggling dc fashion space more than double valuation to b with second raise in month trump urge the uk to drill more oil criticizing it energy policy global startup accelerator techstars officially land in tampa bay bar k bar and dog park combo close door amid financial challenge hospital in region ranked among top in florida north korea demand u recognize it a irreversible nuclear power microsofts msft ai revolution could spark a trilliondollar surge say wedbush mena renewable energy market to reach billion by driving a regional green energy boom bcc research llc mena renewable energy market to reach billion by driving a regional green energy boom bcc research llc stock market today dow fall point sp and nasdaq edge lower ahead of fed meeting and big tech earnings u and china close to extending tariff truce marketwatch brazil harvest pressure and easing tariff fear undercut coffee price google exec say employee have to be more aisavvy a competition ramp up animal shelter are closing their door to surrender a pet care cost rise how owner can afford emergency disney sued for defamation over fx series say nothing yesteryear johnson johnson nysejnj holding lowered by new advisory service llc canopy partner llc ha holding in netflix inc nasdaqnflx top lithium stock to add to your watchlist july th meeder asset management inc grows stock holding in johnson johnson nysejnj warner record name new marketing head lay off in exec realignment visualizing how mass shooting in compare with past year ivory coast president will seek a fourth term after the disqualification of rival why tapestry stock is cooling off despite coach hot streak amazon amzn price target raised to here whats driving it rosen a leading law firm encourages reddit inc investor to secure counsel before important anthropic in talk to raise fresh capital at billion valuation thermo fisher post q rise boeing report smaller loss a ceo see progress in turnaround boeing report smaller loss a ceo see progress in turnaround tax break for charitable contribution are returning for million american tax break for charitable contribution are returning for million american tax break for charitable contribution are returning for million american tax break for charitable contribution are returning for million american boeing share sink despite upbeat ndquarter earnings nick the greek to join burbank business landscape median home price by state how much do house cost rosen a leading law firm encourages reddit inc investor to secure counsel before important deadline in security class action rddt why krispy kreme stock is plummeting today inspector visit lancashire business here how each one scored for food hygiene florida struggle to keep younger worker creating headache for employer lafayette county property transfer honolulu lawsuit against fossil fuel company lead climate change legal fight utah husband wife accused of selling narcotic out of home with child letter local impact pennsylvania legislation aim to protect homebuyers in floodprone neighborhood a qa with historian marco antonio samaniego on the new binational tijuana river agreement dc council launch public hearing on billion rfk stadium proposal letter revelation show treason of those involved in russia hoax nyc subway mystery power outage ha rush hour rider vowing to steal ride a payback bipartisan bill call for annual treasury stress test pete buttigieg say gop epstein scandal reveals breakdown in societal trust trump reacts to starmer plan to recognize palestine could make the case that youre rewarding hamas from the grca can you or your organization afford to ignore professional development gop senator float broken plan to buy off american mad about tariff walking school bus and teaching student their right how educator are preparing for back to school after ice raid riverside revisits nuisance ordinance karrueche tran stayed by deion sander side during lifesaving surgery here more about her trump admin seizes judicial thorn in their side letter decriminalize possessing small amount of cannabis lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity netanyahu is reportedly planning to annex gaza strip with trump admins backing letter america need return to tradition law and order meeder asset management inc grows stake in exxon mobil corporation nysexom top lithium stock to add to your watchlist july th visualizing how mass shooting in compare with past year maxwell lawyer issue list of demand for congressional testimony lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity lawyer for epstein former girlfriend say shes open to interview with congress if given immunity what american think about israel military action in gaza according to a new gallup poll what american think about israel military action in gaza according to a new gallup poll what american think about israel military action in gaza according to a new gallup poll what american think about israel military action in gaza according to a new gallup poll wire fire cause power outage in charlestown and dedham rosen a leading law firm encourages reddit inc investor to secure counsel before important deadline in security class action rddt sydney sweeneys ad spark absurd nazi claim on good morning america legacy medias latest meltdown berger montague pc investigating claim on behalf of sable offshore corp nyse soc investor after class action filing minnesota pause payment to housing stability provider amid fraud investigation what american think about israel military action in gaza according to a new gallup poll what american think about israel military action in gaza according to a new gallup poll what american think about israel military action in gaza according to a new gallup poll israeli ambassador push back on claim that israel is preventing aid distribution in gaza trump proposes to repeal landmark finding that greenhouse gas harm the public trump proposes to repeal landmark finding that greenhouse gas harm the public trump proposes to repeal landmark finding that greenhouse gas harm the public trump proposes to repeal landmark finding that greenhouse gas harm the public trump proposes to repeal landmark finding that greenhouse gas harm the public trump proposes to repeal landmark finding that greenhouse gas harm the public la county move ahead with law requiring ice agent show face badge and agency affiliation but legal hurdle loom house oversight committee reject ghislaine maxwell bid for immunity in exchange for testimony house oversight committee reject ghislaine maxwell request for immunity in exchange for her testimony biggest source of immigrant to new hampshire biggest source of immigrant to kansa biggest source of immigrant to kansa harvest fund management co ltd raise stake in alphabet inc nasdaqgoog biggest source of immigrant to alabama biggest source of immigrant to maine ucla to pay million to settle complaint of discrimination against jewish student ucla to pay million to settle complaint of discrimination against jewish student dont question my integrity spartacus return in dramatic senate floor showdown for jersey a daca recipient accidentally drove into mexico now he being fasttracked for deportation cashflow underwriting would bring fairness to mortgage lending the independent view big beautiful bill a former cdc vaccine adviser ha word for rfk jr segwayriding pest control salesperson described a pest themselves in upper township elsewhere gov walz activates minnesota national guard over st paul cyberattack gov walz activates minnesota national guard over st paul cyberattack nantucket official accuse offshore wind developer of going into hiding since trump election nantucket official accuse offshore wind developer of going into hiding since trump election nantucket official accuse offshore wind developer of going into hiding since trump election nc u senate race set to be one of nation biggest wrong way tom zirpoli big government is back commentary guest columnist an opportunity for much needed clarification together review the only one for me is you and you for me here how to see this week double meteor shower man run for mayor with really really really unfortunate name justice department fire two senior antitrust attorney alleging insubordination democraticled state sue trump admin over law that defunds planned parenthood uk will recognize palestine a a state if ceasefire not agreed on pm say mackenzie financial corp sell share of intuitive surgical inc nasdaqisrg trump say india will likely face a to tariff yearold girl reported missing in worcester'
$ws.Cells.Item(26, 11).Value = 0.115
$ws.Cells.Item(26, 12).Value = 0.096
$ws.Cells.Item(26, 13).Value = 0.789
$ws.Cells.Item(26, 14).Value = 0.9888
$ws.Cells.Item(26, 15).Value = 141
$ws.Cells.Item(26, 16).Value = 95
$ws.Cells.Item(26, 17).Value = 0
$ws.Cells.Item(26, 18).Value = 163
$ws.Cells.Item(26, 19).Value = 55
$ws.Cells.Item(26, 20).Value = 175
$ws.Cells.Item(26, 21).Value = 157
$ws.Cells.Item(26, 22).Value = 58
$ws.Cells.Item(26, 23).Value = 82
$ws.Cells.Item(26, 24).Value = 37

# --- Row 27 (2025-07-31) ---
$ws.Cells.Item(27, 1).Value = 45869
$ws.Cells.Item(27, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(27, 2).Value = 6339.39013671875
$ws.Cells.Item(27, 3).Value = 6427.02001953125
$ws.Cells.Item(27, 4).Value = 6327.64013671875
$ws.Cells.Item(27, 5).Value = 6427.02001953125
$ws.Cells.Item(27, 6).Value = 4035281000
$ws.Cells.Item(27, 7).Value = -0.003694819341153
$ws.Cells.Item(27, 8).Value = -1
$ws.Cells.Item(27, 9).Value = 'boeing narrow loss ahead of new cfo arrival how cadiz plan to deliver water relief to droughtprone community nvidias nvda ai chip in high demand despite export limit jefferies weighs in gaza reaching tipping point with clinic full of malnourished child say save the child worker we are doing it with a heavy heart israeli right group accuse government of genocide ap news summary at pm edt hoffa legend endures year after exdetroit union leader disappearance shiba inu news today shib price face key resistance a trader await august breakout chatgpt ceo sends a stern warning on this alarming practice meta q earnings update wall street is bullish eyeing ai opportunity with stock up this year business insider executive contract whats common and what to look for man in custody after a virginia council member wa set on fire ncdhhs partner with talkspace to provide free virtual mental health therapy to teen involved in or impacted by the justice system wall street stick with apple aapl a earnings and ai question loom ap technology summarybrief at pm edt edge fall protection expands headquarters and team to support continued growth fed hold main rate steady overriding two dissenting vote see the global impact of the sixthlargest earthquake ever ruvi ai ruvi just caught whale attention coinmarketcap listing spark frenzy could turn into figure unlike cardano ada inventhelp inventor develops new warming cooling vest ctk best coin under ozak ai could jump in one month new mexicobased wesst to close office cut staff amid federal funding uncertainty amazon sleek and stylish vintage laptop backpack is on sale for and shopper say it fit everything best portable bluetooth speaker of cpkc see profit and revenue jump on stronger volume barclays warns trader against chasing high oil price providence st mary foundation golf tournament raise adcolor announces nominee for th annual adcolor award fm to conclude with expertled ask anything closing session what are the advantage and disadvantage of artificial intelligence rebuilding america how skilled trade can revive community and career alibaba baba enters the ai wearable race with new smart glass divided fed hold key interest rate steady defying trump demand for aggressive cut whats going on with meta stock wednesday stronger sugar output in brazil pressure sugar price elon musk plan to build boring co tunnel in nashville spark partisan feud gemini cli hack expose critical security flaw in coding tool oklo stock is surging wednesday whats driving the action your labor day shopping guide to the best endofsummer deal high noon say some vodka seltzer were mislabeled a celsius energy drink issue recall high noon recall vodka seltzer pack mislabeled a celsius energy drink cbs news guest column minimum wage falling further behind cost of living solid july job gain reinforce fed wait and see rate stance how deeb ghezlan is expanding capital access for america real estate operator ap business summarybrief at pm edt dollar strength undercut ny cocoa price dont steal that hotel pillow buy it here for off watch live too late powell explains why ratecuts remain off the table despite dissent growth doubt federal reserve leaf interest rate unchanged even a trump demand cut how stablecoins are transforming enterprise financial strategy the fed keep borrowing cost high despite trump push for cutsheres why amazon to rally more than here are top analyst forecast for wednesday nba veteran charged with defrauding over k from strip casino ge vernova gev get mixed analyst call downgrade with a higher price target fed skip rate cut againbut reveals something we havent seen since reduced tariff risk weigh on arabica coffee white house trump suspends de minimis exemption for commercial shipment globally construction pro react to interest rate decision trump announces tariff on india and unspecified penalty for buying russian oil cnbc host mercilessly mock elizabeth warren and dems economic fearmongering v reality whats going on with coinbase stock today clean energy tax credit are going away here how to get them before it too late automaker are eating the cost of tariff for now npr what is a home equity agreement what is a home equity agreement what is a home equity agreement aw shuck national oyster day is tuesday aug who owns authorized user reward sometimes it complicated u senator raise alarm over delta air line testing of personalized pricing servicenow now downgraded due to limited upside potential hawthorn bancshares report second quarter result walmart flash deal give massive saving on kitchen essential this week army secretary withdraws west point job offer to former biden official amid pressure from farright activist hawthorn bancshares report second quarter result cnbc daily open the tariff rate you get is not just about your trade policy citi on board a transfer partner with american airline u stock market hold steady after the fed make no move on interest rate u stock market hold steady after the fed make no move on interest rate u stock market hold steady after the fed make no move on interest rate u stock market hold steady after the fed make no move on interest rate u stock market hold steady after the fed make no move on interest rate u stock market hold steady after the fed make no move on interest rate worker sue hormel alleging violation of minnesota new sick and safe time law federal reserve vote on interest rate illinois tool work achieves record q boost outlook tesla tsla get analyst support a robotaxi vision move toward bay area launch ebay stock surge on earnings beat rosy guidance tc layoff expose deeper workforce crisis than rival analyst weigh in invesco canada announces risk rating change to canadian mutual fund fed decides against interest rate cutbut powell face first double dissent in decade u economy slowed in first half of a tariff scrambled data the new york time senate help committee unanimously pass esop association legislative priority the fed stay on hold and the key number to watch when meta and microsoft report cnbc the fed stay on hold and the key number to watch when meta and microsoft report federal reserve hold it benchmark rate steady at today fomc meeting fed leaf fund rate unchanged ap news summary at pm edt federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut some small maine college football program can finally join the playoff party it about time foreign couple flock to denmark to get married copenhagen want to save room for local foreign couple flock to denmark to get married copenhagen want to save room for local foreign couple flock to denmark to get married copenhagen want to save room for local foreign couple flock to denmark to get married copenhagen want to save room for local foreign couple flock to denmark to get married copenhagen want to save room for local gaza reaching tipping point with clinic full of malnourished child say save the child worker we are doing it with a heavy heart israeli right group accuse government of genocide newport beach fiscal folly a it pay two city manager live now senate foreign relation subcommittee hearing on finding opportunity for usafrica partnership man in custody after a virginia council member wa set on fire see the global impact of the sixthlargest earthquake ever wa review reconnecting vancouver wildlife area with columbia river whats the difference between shame and guilt the ethical life podcast whats the difference between shame and guilt the ethical life podcast whats the difference between shame and guilt the ethical life podcast illinois lawmaker official seek input into how auto insurance rate are set whats the difference between shame and guilt the ethical life podcast whats the difference between shame and guilt the ethical life podcast whats the difference between shame and guilt the ethical life podcast sc say lack of scientific standard prevents testing for freshwater amoeba that attack brain whats the difference between shame and guilt the ethical life podcast pope mark th anniversary of helsinki accord school and county begin to see payment delay a pennsylvania budget stalemate hit a month democrat blast politics of fear after trump claim against obama top aide a tiny california city prone to tsunami saw the highest wave in continental u after quake he always too late trump continues criticism of powell over interest rate do to require inperson consular interview for most nonimmigrant visa applicant trump administration is launching a new private health tracking system with big tech help trump unveils south korea trade deal boasting b investment in u gop sen josh hawley democrat vote to advance congressional stock trading ban texas republican unveil congressional map giving them a chance to pick up seat trump announces tariff on india and unspecified penalty for buying russian oil loyalty over democracy why is the democratic party favorability rating so low senate committee advance trump nominee to lead cybersecurity agency that protects election system virginia city councilman set on fire in local magazine office virginia city councilman set on fire in workplace attack suspect charged with attempted murder automaker are eating the cost of tariff for now kamala harris decides against run for california governor in doj say releasing epstein transcript serf public interest in new court filing federal reserve keep interest rate unchanged for fifth time army secretary withdraws west point job offer to former biden official amid pressure from farright activist analysisus reversal on key climate finding spell uncertainty for business u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro fed decides against interest rate cutbut powell face first double dissent in decade exclusive trump issue veto threat against bernie sander latest antiisrael move senate help committee unanimously pass esop association legislative priority federal reserve leaf interest rate unchanged even a trump demand cut state rep face new shoplifting charge senate panel advance bill banning congressional stock trading why harrison ford think politics may be making a healthy swing to the right scammer play into confusion over tariff immigration drug trafficking fed hold rate steady but see first double dissent in three decade gloversville republican committee endorses antonucci for mayor fed hold rate steady but see first double dissent in three decade fed hold rate steady but see first double dissent in three decade fed hold rate steady but see first double dissent in three decade fed hold rate steady but see first double dissent in three decade fed hold rate steady but see first double dissent in three decade whither the democratic party editorial gop senate leader reject trump call for end to bipartisan judicial nomination tradition epstein and maxwell grand jury transcript include testimony from only two witness justice department say texas republican propose new u house map with more winnable gop seat'
$ws.Cells.Item(27, 10).Value = 'boeing narrow loss ahead of new cfo arrival how cadiz plan to deliver water relief to droughtprone community nvidias nvda ai chip in high demand despite export limit jefferies weighs in gaza reaching tipping point with clinic full of malnourished child say save the child worker we are doing it with a heavy heart israeli right group accuse government of genocide ap news summary at pm edt hoffa legend endures year after exdetroit union leader disappearance shiba inu news today shib price face key resistance a trader await august breakout chatgpt ceo sends a stern warning on this alarming practice meta q earnings update wall street is bullish eyeing ai opportunity with stock up this year business insider executive contract whats common and what to look for man in custody after a virginia council member wa set on fire ncdhhs partner with talkspace to provide free virtual mental health therapy to teen involved in or impacted by the justice system wall street stick with apple aapl a earnings and ai question loom ap technology summarybrief at pm edt edge fall protection expands headquarters and team to support continued growth fed hold main rate steady overriding two dissenting vote see the global impact of the sixthlargest earthquake ever ruvi ai ruvi just caught whale attention coinmarketcap listing spark frenzy could turn into figure unlike cardano ada inventhelp inventor develops new warming cooling vest ctk best coin under ozak ai could jump in one month new mexicobased wesst to close office cut staff amid federal funding uncertainty amazon sleek and stylish vintage laptop backpack is on sale for and shopper say it fit everything best portable bluetooth speaker of cpkc see profit and revenue jump on stronger volume barclays warns trader against chasing high oil price providence st mary foundation golf tournament raise adcolor announces nominee for th annual adcolor award fm to conclude with expertled ask anything closing session what are the advantage and disadvantage of artificial intelligence rebuilding america how skilled trade can revive community and career alibaba baba enters the ai wearable race with new smart glass divided fed hold key interest rate steady defying trump demand for aggressive cut whats going on with meta stock wednesday stronger sugar output in brazil pressure sugar price elon musk plan to build boring co tunnel in nashville spark partisan feud gemini cli hack expose critical security flaw in coding tool oklo stock is surging wednesday whats driving the action your labor day shopping guide to the best endofsummer deal high noon say some vodka seltzer were mislabeled a celsius energy drink issue recall high noon recall vodka seltzer pack mislabeled a celsius energy drink cbs news guest column minimum wage falling further behind cost of living solid july job gain reinforce fed wait and see rate stance how deeb ghezlan is expanding capital access for america real estate operator ap business summarybrief at pm edt dollar strength undercut ny cocoa price dont steal that hotel pillow buy it here for off watch live too late powell explains why ratecuts remain off the table despite dissent growth doubt federal reserve leaf interest rate unchanged even a trump demand cut how stablecoins are transforming enterprise financial strategy the fed keep borrowing cost high despite trump push for cutsheres why amazon to rally more than here are top analyst forecast for wednesday nba veteran charged with defrauding over k from strip casino ge vernova gev get mixed analyst call downgrade with a higher price target fed skip rate cut againbut reveals something we havent seen since reduced tariff risk weigh on arabica coffee white house trump suspends de minimis exemption for commercial shipment globally construction pro react to interest rate decision trump announces tariff on india and unspecified penalty for buying russian oil cnbc host mercilessly mock elizabeth warren and dems economic fearmongering v reality whats going on with coinbase stock today clean energy tax credit are going away here how to get them before it too late automaker are eating the cost of tariff for now npr what is a home equity agreement what is a home equity agreement what is a home equity agreement aw shuck national oyster day is tuesday aug who owns authorized user reward sometimes it complicated u senator raise alarm over delta air line testing of personalized pricing servicenow now downgraded due to limited upside potential hawthorn bancshares report second quarter result walmart flash deal give massive saving on kitchen essential this week army secretary withdraws west point job offer to former biden official amid pressure from farright activist hawthorn bancshares report second quarter result cnbc daily open the tariff rate you get is not just about your trade policy citi on board a transfer partner with american airline u stock market hold steady after the fed make no move on interest rate u stock market hold steady after the fed make no move on interest rate u stock market hold steady after the fed make no move on interest rate u stock market hold steady after the fed make no move on interest rate u stock market hold steady after the fed make no move on interest rate u stock market hold steady after the fed make no move on interest rate worker sue hormel alleging violation of minnesota new sick and safe time law federal reserve vote on interest rate illinois tool work achieves record q boost outlook tesla tsla get analyst support a robotaxi vision move toward bay area launch ebay stock surge on earnings beat rosy guidance tc layoff expose deeper workforce crisis than rival analyst weigh in invesco canada announces risk rating change to canadian mutual fund fed decides against interest rate cutbut powell face first double dissent in decade u economy slowed in first half of a tariff scrambled data the new york time senate help committee unanimously pass esop association legislative priority the fed stay on hold and the key number to watch when meta and microsoft report cnbc the fed stay on hold and the key number to watch when meta and microsoft report federal reserve hold it benchmark rate steady at today fomc meeting fed leaf fund rate unchanged ap news summary at pm edt federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut federal reserve leaf interest rate unchanged even a trump demand cut some small maine college football program can finally join the playoff party it about time foreign couple flock to denmark to get married copenhagen want to save room for local foreign couple flock to denmark to get married copenhagen want to save room for local foreign couple flock to denmark to get married copenhagen want to save room for local foreign couple flock to denmark to get married copenhagen want to save room for local foreign couple flock to denmark to get married copenhagen want to save room for local gaza reaching tipping point with clinic full of malnourished child say save the child worker we are doing it with a heavy heart israeli right group accuse government of genocide newport beach fiscal folly a it pay two city manager live now senate foreign relation subcommittee hearing on finding opportunity for usafrica partnership man in custody after a virginia council member wa set on fire see the global impact of the sixthlargest earthquake ever wa review reconnecting vancouver wildlife area with columbia river whats the difference between shame and guilt the ethical life podcast whats the difference between shame and guilt the ethical life podcast whats the difference between shame and guilt the ethical life podcast illinois lawmaker official seek input into how auto insurance rate are set whats the difference between shame and guilt the ethical life podcast whats the difference between shame and guilt the ethical life podcast whats the difference between shame and guilt the ethical life podcast sc say lack of scientific standard prevents testing for freshwater amoeba that attack brain whats the difference between shame and guilt the ethical life podcast pope mark th anniversary of helsinki accord school and county begin to see payment delay a pennsylvania budget stalemate hit a month democrat blast politics of fear after trump claim against obama top aide a tiny california city prone to tsunami saw the highest wave in continental u after quake he always too late trump continues criticism of powell over interest rate do to require inperson consular interview for most nonimmigrant visa applicant trump administration is launching a new private health tracking system with big tech help trump unveils south korea trade deal boasting b investment in u gop sen josh hawley democrat vote to advance congressional stock trading ban texas republican unveil congressional map giving them a chance to pick up seat trump announces tariff on india and unspecified penalty for buying russian oil loyalty over democracy why is the democratic party favorability rating so low senate committee advance trump nominee to lead cybersecurity agency that protects election system virginia city councilman set on fire in local magazine office virginia city councilman set on fire in workplace attack suspect charged with attempted murder automaker are eating the cost of tariff for now kamala harris decides against run for california governor in doj say releasing epstein transcript serf public interest in new court filing federal reserve keep interest rate unchanged for fifth time army secretary withdraws west point job offer to former biden official amid pressure from farright activist analysisus reversal on key climate finding spell uncertainty for business u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro u sanction brazil supreme court justice overseeing case against bolsonaro fed decides against interest rate cutbut powell face first double dissent in decade exclusive trump issue veto threat against bernie sander latest antiisrael move senate help committee unanimously pass esop association legislative priority federal reserve leaf interest rate unchanged even a trump demand cut state rep face new shoplifting charge senate panel advance bill banning congressional stock trading why harrison ford think politics may be making a healthy swing to the right scammer play into confusion over tariff immigration drug trafficking fed hold rate steady but see first double dissent in three decade gloversville republican committee endorses antonucci for mayor fed hold rate steady but see first double dissent in three decade fed hold rate steady but see first double dissent in three decade fed hold rate steady but see first double dissent in three decade fed hold rate steady but see first double dissent in three decade fed hold rate steady but see first double dissent in three decade whither the democratic party editorial gop senate leader reject trump call for end to bipartisan judicial nomination tradition epstein and maxwell grand jury transcript include testimony from only two witness justice department say texas republican propose new u house map with more winnable gop seat'
$ws.Cells.Item(27, 11).Value = 0.142
$ws.Cells.Item(27, 12).Value = 0.082
$ws.Cells.Item(27, 13).Value = 0.776
$ws.Cells.Item(27, 14).Value = 0.9992
$ws.Cells.Item(27, 15).Value = 95
$ws.Cells.Item(27, 16).Value = 88
$ws.Cells.Item(27, 17).Value = 0
$ws.Cells.Item(27, 18).Value = 125
$ws.Cells.Item(27, 19).Value = 74
$ws.Cells.Item(27, 20).Value = 228
$ws.Cells.Item(27, 21).Value = 143
$ws.Cells.Item(27, 22).Value = 62
$ws.Cells.Item(27, 23).Value = 29
$ws.Cells.Item(27, 24).Value = 36

